# TAWA_TestCases.xlsx - "test case is updated"
#
# Adds a new "Test Data" column (between the existing "Precondition" and
# "Steps" columns) to the ReservedTrips and Booking sheets, and updates the
# active sheet/selection so ReservedTrips (rather than Login) is the tab
# that is selected when the workbook is reopened.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Booking sheet: insert the new "Test Data" column (E) and move the
# selection to F3.
# ---------------------------------------------------------------------
$wsBooking = $wb.Worksheets.Item("Booking")

$wsBooking.Columns("E").Insert()
$wsBooking.Range("E1").Value = "Test Data"
$wsBooking.Columns("E").ColumnWidth = $wsBooking.Columns("D").ColumnWidth

$wsBooking.Range("F3").Select()

# ---------------------------------------------------------------------
# ReservedTrips sheet: insert the new "Test Data" column (E), move the
# selection to H5, and make this the active/selected sheet (selecting a
# range on it last is what flips tabSelected / the workbook's active tab).
# ---------------------------------------------------------------------
$wsReserved = $wb.Worksheets.Item("ReservedTrips")

$wsReserved.Columns("E").Insert()
$wsReserved.Range("E1").Value = "Test Data"
$wsReserved.Columns("E").ColumnWidth = $wsReserved.Columns("D").ColumnWidth

$wsReserved.Activate()
$wsReserved.Range("H5").Select()
